{"js": "// Auto-generated edit: replace each old math-fact / date string with its\n// corresponding new value. Each \"old\" text is unique in the document, so a\n// targeted body.search(...) + insertText(..., \"Replace\") pair is safe and\n// avoids disturbing any other run formatting.\nconst replacements = [\n  [\n    \"2024-04-14 Sunday\",\n    \"2024-04-15 Monday\"\n  ],\n  [\n    \"67-29=38\",\n    \"55+3=58\"\n  ],\n  [\n    \"49+27=76\",\n    \"10+9=19\"\n  ],\n  [\n    \"49+44=93\",\n    \"86-25=61\"\n  ],\n  [\n    \"8-8=0\",\n    \"69-44=25\"\n  ],\n  [\n    \"68-24=44\",\n    \"43+15=58\"\n  ],\n  [\n    \"72-37=35\",\n    \"80-4=76\"\n  ],\n  [\n    \"74-40=34\",\n    \"17+12=29\"\n  ],\n  [\n    \"48-37=11\",\n    \"78+5=83\"\n  ],\n  [\n    \"7+67=74\",\n    \"92-70=22\"\n  ],\n  [\n    \"48-10=38\",\n    \"74+18=92\"\n  ],\n  [\n    \"30+2=32\",\n    \"52-44=8\"\n  ],\n  [\n    \"99-2=97\",\n    \"17+13=30\"\n  ],\n  [\n    \"75-0=75\",\n    \"30+26=56\"\n  ],\n  [\n    \"23+57=80\",\n    \"37+29=66\"\n  ],\n  [\n    \"43-9=34\",\n    \"52-10=42\"\n  ],\n  [\n    \"52+34=86\",\n    \"81-11=70\"\n  ],\n  [\n    \"84-13=71\",\n    \"7+23=30\"\n  ],\n  [\n    \"66+0=66\",\n    \"6-0=6\"\n  ],\n  [\n    \"45-12=33\",\n    \"36+8=44\"\n  ],\n  [\n    \"4+35=39\",\n    \"30+26=56\"\n  ],\n  [\n    \"58+14=72\",\n    \"61-10=51\"\n  ],\n  [\n    \"53-30=23\",\n    \"48+49=97\"\n  ],\n  [\n    \"37+7=44\",\n    \"39-11=28\"\n  ],\n  [\n    \"97-60=37\",\n    \"44+47=91\"\n  ],\n  [\n    \"66+30=96\",\n    \"81-22=59\"\n  ],\n  [\n    \"69-68=1\",\n    \"61+4=65\"\n  ],\n  [\n    \"92-45=47\",\n    \"14+56=70\"\n  ],\n  [\n    \"40+35=75\",\n    \"38-0=38\"\n  ],\n  [\n    \"60-53=7\",\n    \"55-34=21\"\n  ],\n  [\n    \"61-13=48\",\n    \"63+21=84\"\n  ],\n  [\n    \"56+19=75\",\n    \"52-46=6\"\n  ],\n  [\n    \"49-37=12\",\n    \"48+51=99\"\n  ],\n  [\n    \"17+70=87\",\n    \"11+76=87\"\n  ],\n  [\n    \"16+24=40\",\n    \"11+16=27\"\n  ],\n  [\n    \"87-48=39\",\n    \"17+68=85\"\n  ],\n  [\n    \"1+66=67\",\n    \"8+48=56\"\n  ],\n  [\n    \"75-43=32\",\n    \"2+49=51\"\n  ],\n  [\n    \"8+31=39\",\n    \"27+8=35\"\n  ],\n  [\n    \"70-36=34\",\n    \"98-75=23\"\n  ],\n  [\n    \"55-25=30\",\n    \"96-0=96\"\n  ],\n  [\n    \"67+25=92\",\n    \"40+19=59\"\n  ],\n  [\n    \"76-51=25\",\n    \"93-71=22\"\n  ],\n  [\n    \"55+33=88\",\n    \"88-69=19\"\n  ],\n  [\n    \"40-25=15\",\n    \"48-32=16\"\n  ],\n  [\n    \"6+27=33\",\n    \"58-28=30\"\n  ],\n  [\n    \"35+39=74\",\n    \"2+8=10\"\n  ],\n  [\n    \"57-5=52\",\n    \"50-26=24\"\n  ],\n  [\n    \"52-36=16\",\n    \"76-72=4\"\n  ],\n  [\n    \"88-34=54\",\n    \"20+35=55\"\n  ],\n  [\n    \"24-13=11\",\n    \"92-37=55\"\n  ],\n  [\n    \"82-55=27\",\n    \"43-4=39\"\n  ],\n  [\n    \"5+38=43\",\n    \"0+99=99\"\n  ],\n  [\n    \"9+83=92\",\n    \"27-19=8\"\n  ],\n  [\n    \"36+4=40\",\n    \"79-47=32\"\n  ],\n  [\n    \"7+4=11\",\n    \"34+25=59\"\n  ],\n  [\n    \"92-49=43\",\n    \"58+4=62\"\n  ],\n  [\n    \"98-27=71\",\n    \"62-0=62\"\n  ],\n  [\n    \"26+69=95\",\n    \"37+45=82\"\n  ],\n  [\n    \"52-40=12\",\n    \"63-19=44\"\n  ],\n  [\n    \"4+34=38\",\n    \"1+30=31\"\n  ],\n  [\n    \"14+70=84\",\n    \"71+18=89\"\n  ],\n  [\n    \"6+77=83\",\n    \"38-32=6\"\n  ],\n  [\n    \"78-13=65\",\n    \"62+25=87\"\n  ],\n  [\n    \"48-31=17\",\n    \"80-46=34\"\n  ],\n  [\n    \"19+71=90\",\n    \"66-29=37\"\n  ],\n  [\n    \"24+60=84\",\n    \"9+37=46\"\n  ],\n  [\n    \"36-30=6\",\n    \"31-31=0\"\n  ],\n  [\n    \"59-39=20\",\n    \"39+48=87\"\n  ],\n  [\n    \"88+8=96\",\n    \"47+2=49\"\n  ],\n  [\n    \"88-50=38\",\n    \"37-27=10\"\n  ],\n  [\n    \"71-60=11\",\n    \"9+19=28\"\n  ],\n  [\n    \"41+56=97\",\n    \"71-64=7\"\n  ],\n  [\n    \"71-49=22\",\n    \"8+47=55\"\n  ],\n  [\n    \"39+53=92\",\n    \"0+73=73\"\n  ],\n  [\n    \"53-6=47\",\n    \"89-51=38\"\n  ],\n  [\n    \"66-5=61\",\n    \"40+37=77\"\n  ],\n  [\n    \"31+37=68\",\n    \"97-66=31\"\n  ],\n  [\n    \"3+50=53\",\n    \"72-0=72\"\n  ],\n  [\n    \"86-26=60\",\n    \"48-8=40\"\n  ],\n  [\n    \"20+6=26\",\n    \"75-34=41\"\n  ],\n  [\n    \"40-29=11\",\n    \"24-19=5\"\n  ],\n  [\n    \"47-5=42\",\n    \"8+89=97\"\n  ],\n  [\n    \"13-10=3\",\n    \"39+28=67\"\n  ],\n  [\n    \"81+9=90\",\n    \"36+21=57\"\n  ],\n  [\n    \"52-7=45\",\n    \"71-4=67\"\n  ],\n  [\n    \"18+29=47\",\n    \"72-18=54\"\n  ],\n  [\n    \"27+33=60\",\n    \"83-0=83\"\n  ],\n  [\n    \"22-18=4\",\n    \"86-49=37\"\n  ],\n  [\n    \"3+35=38\",\n    \"1+24=25\"\n  ],\n  [\n    \"34+27=61\",\n    \"7+9=16\"\n  ],\n  [\n    \"5+73=78\",\n    \"26+55=81\"\n  ],\n  [\n    \"74+16=90\",\n    \"46-43=3\"\n  ],\n  [\n    \"24+32=56\",\n    \"89-46=43\"\n  ],\n  [\n    \"27+22=49\",\n    \"57+24=81\"\n  ],\n  [\n    \"95-75=20\",\n    \"64+11=75\"\n  ],\n  [\n    \"23+73=96\",\n    \"89-10=79\"\n  ],\n  [\n    \"82-54=28\",\n    \"41-23=18\"\n  ],\n  [\n    \"6+68=74\",\n    \"83-28=55\"\n  ],\n  [\n    \"80-78=2\",\n    \"28+28=56\"\n  ],\n  [\n    \"72-10=62\",\n    \"95-57=38\"\n  ]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Auto-generated edit: replace each old math-fact / date string with its\n# corresponding new value using Word's Find & Replace (wdReplaceOne), one\n# pair at a time so that each replacement targets exactly the single\n# occurrence of its (unique) old text.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-04-14 Sunday\", \"2024-04-15 Monday\"),\n    @(\"67-29=38\", \"55+3=58\"),\n    @(\"49+27=76\", \"10+9=19\"),\n    @(\"49+44=93\", \"86-25=61\"),\n    @(\"8-8=0\", \"69-44=25\"),\n    @(\"68-24=44\", \"43+15=58\"),\n    @(\"72-37=35\", \"80-4=76\"),\n    @(\"74-40=34\", \"17+12=29\"),\n    @(\"48-37=11\", \"78+5=83\"),\n    @(\"7+67=74\", \"92-70=22\"),\n    @(\"48-10=38\", \"74+18=92\"),\n    @(\"30+2=32\", \"52-44=8\"),\n    @(\"99-2=97\", \"17+13=30\"),\n    @(\"75-0=75\", \"30+26=56\"),\n    @(\"23+57=80\", \"37+29=66\"),\n    @(\"43-9=34\", \"52-10=42\"),\n    @(\"52+34=86\", \"81-11=70\"),\n    @(\"84-13=71\", \"7+23=30\"),\n    @(\"66+0=66\", \"6-0=6\"),\n    @(\"45-12=33\", \"36+8=44\"),\n    @(\"4+35=39\", \"30+26=56\"),\n    @(\"58+14=72\", \"61-10=51\"),\n    @(\"53-30=23\", \"48+49=97\"),\n    @(\"37+7=44\", \"39-11=28\"),\n    @(\"97-60=37\", \"44+47=91\"),\n    @(\"66+30=96\", \"81-22=59\"),\n    @(\"69-68=1\", \"61+4=65\"),\n    @(\"92-45=47\", \"14+56=70\"),\n    @(\"40+35=75\", \"38-0=38\"),\n    @(\"60-53=7\", \"55-34=21\"),\n    @(\"61-13=48\", \"63+21=84\"),\n    @(\"56+19=75\", \"52-46=6\"),\n    @(\"49-37=12\", \"48+51=99\"),\n    @(\"17+70=87\", \"11+76=87\"),\n    @(\"16+24=40\", \"11+16=27\"),\n    @(\"87-48=39\", \"17+68=85\"),\n    @(\"1+66=67\", \"8+48=56\"),\n    @(\"75-43=32\", \"2+49=51\"),\n    @(\"8+31=39\", \"27+8=35\"),\n    @(\"70-36=34\", \"98-75=23\"),\n    @(\"55-25=30\", \"96-0=96\"),\n    @(\"67+25=92\", \"40+19=59\"),\n    @(\"76-51=25\", \"93-71=22\"),\n    @(\"55+33=88\", \"88-69=19\"),\n    @(\"40-25=15\", \"48-32=16\"),\n    @(\"6+27=33\", \"58-28=30\"),\n    @(\"35+39=74\", \"2+8=10\"),\n    @(\"57-5=52\", \"50-26=24\"),\n    @(\"52-36=16\", \"76-72=4\"),\n    @(\"88-34=54\", \"20+35=55\"),\n    @(\"24-13=11\", \"92-37=55\"),\n    @(\"82-55=27\", \"43-4=39\"),\n    @(\"5+38=43\", \"0+99=99\"),\n    @(\"9+83=92\", \"27-19=8\"),\n    @(\"36+4=40\", \"79-47=32\"),\n    @(\"7+4=11\", \"34+25=59\"),\n    @(\"92-49=43\", \"58+4=62\"),\n    @(\"98-27=71\", \"62-0=62\"),\n    @(\"26+69=95\", \"37+45=82\"),\n    @(\"52-40=12\", \"63-19=44\"),\n    @(\"4+34=38\", \"1+30=31\"),\n    @(\"14+70=84\", \"71+18=89\"),\n    @(\"6+77=83\", \"38-32=6\"),\n    @(\"78-13=65\", \"62+25=87\"),\n    @(\"48-31=17\", \"80-46=34\"),\n    @(\"19+71=90\", \"66-29=37\"),\n    @(\"24+60=84\", \"9+37=46\"),\n    @(\"36-30=6\", \"31-31=0\"),\n    @(\"59-39=20\", \"39+48=87\"),\n    @(\"88+8=96\", \"47+2=49\"),\n    @(\"88-50=38\", \"37-27=10\"),\n    @(\"71-60=11\", \"9+19=28\"),\n    @(\"41+56=97\", \"71-64=7\"),\n    @(\"71-49=22\", \"8+47=55\"),\n    @(\"39+53=92\", \"0+73=73\"),\n    @(\"53-6=47\", \"89-51=38\"),\n    @(\"66-5=61\", \"40+37=77\"),\n    @(\"31+37=68\", \"97-66=31\"),\n    @(\"3+50=53\", \"72-0=72\"),\n    @(\"86-26=60\", \"48-8=40\"),\n    @(\"20+6=26\", \"75-34=41\"),\n    @(\"40-29=11\", \"24-19=5\"),\n    @(\"47-5=42\", \"8+89=97\"),\n    @(\"13-10=3\", \"39+28=67\"),\n    @(\"81+9=90\", \"36+21=57\"),\n    @(\"52-7=45\", \"71-4=67\"),\n    @(\"18+29=47\", \"72-18=54\"),\n    @(\"27+33=60\", \"83-0=83\"),\n    @(\"22-18=4\", \"86-49=37\"),\n    @(\"3+35=38\", \"1+24=25\"),\n    @(\"34+27=61\", \"7+9=16\"),\n    @(\"5+73=78\", \"26+55=81\"),\n    @(\"74+16=90\", \"46-43=3\"),\n    @(\"24+32=56\", \"89-46=43\"),\n    @(\"27+22=49\", \"57+24=81\"),\n    @(\"95-75=20\", \"64+11=75\"),\n    @(\"23+73=96\", \"89-10=79\"),\n    @(\"82-54=28\", \"41-23=18\"),\n    @(\"6+68=74\", \"83-28=55\"),\n    @(\"80-78=2\", \"28+28=56\"),\n    @(\"72-10=62\", \"95-57=38\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        Write-Output \"WARNING: not found -> $oldText\"\n    }\n}\n\nWrite-Output \"done\"\n"}
